$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 39
$ws.Cells.Item(7, 6).Value = 20
$ws.Cells.Item(8, 6).Value = 2685
$ws.Cells.Item(10, 6).Value = 6362
$ws.Cells.Item(11, 6).Value = 2392
$ws.Cells.Item(13, 6).Value = 22
$ws.Cells.Item(15, 6).Value = 2566
$ws.Cells.Item(16, 6).Value = 29
$ws.Cells.Item(17, 6).Value = 24
$ws.Cells.Item(18, 6).Value = 6772
$ws.Cells.Item(19, 6).Value = 246
$ws.Cells.Item(21, 6).Value = 174
$ws.Cells.Item(22, 6).Value = 108
$ws.Cells.Item(24, 6).Value = 7604
$ws.Cells.Item(25, 6).Value = 13
$ws.Cells.Item(27, 6).Value = 252
$ws.Cells.Item(28, 6).Value = 51
$ws.Cells.Item(32, 6).Value = 46
$ws.Cells.Item(33, 6).Value = 29
$ws.Cells.Item(35, 6).Value = 29
$ws.Cells.Item(37, 6).Value = 63
$ws.Cells.Item(38, 6).Value = 2571
$ws.Cells.Item(39, 6).Value = 37
$ws.Cells.Item(40, 6).Value = 79
$ws.Cells.Item(41, 6).Value = 27
$ws.Cells.Item(42, 6).Value = 1145
$ws.Cells.Item(44, 6).Value = 590
$ws.Cells.Item(45, 6).Value = 3595
$ws.Cells.Item(46, 6).Value = 127
$ws.Cells.Item(47, 6).Value = 1154
$ws.Cells.Item(48, 6).Value = 88

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 24
$ws.Cells.Item(5, 6).Value = 231
$ws.Cells.Item(13, 6).Value = 1

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 39
$ws.Cells.Item(5, 6).Value = 20
$ws.Cells.Item(7, 6).Value = 2685
$ws.Cells.Item(8, 6).Value = 24
$ws.Cells.Item(9, 6).Value = 231
$ws.Cells.Item(10, 6).Value = 6362
$ws.Cells.Item(11, 6).Value = 2392
$ws.Cells.Item(13, 6).Value = 22
$ws.Cells.Item(15, 6).Value = 2566
$ws.Cells.Item(16, 6).Value = 29
$ws.Cells.Item(19, 6).Value = 24
$ws.Cells.Item(20, 6).Value = 6772
$ws.Cells.Item(21, 6).Value = 246
$ws.Cells.Item(23, 6).Value = 174
$ws.Cells.Item(25, 6).Value = 7604
$ws.Cells.Item(26, 6).Value = 13
$ws.Cells.Item(28, 6).Value = 252
$ws.Cells.Item(29, 6).Value = 51
$ws.Cells.Item(34, 6).Value = 29
$ws.Cells.Item(37, 6).Value = 63
$ws.Cells.Item(38, 6).Value = 2571
$ws.Cells.Item(39, 6).Value = 37
$ws.Cells.Item(40, 6).Value = 27
$ws.Cells.Item(41, 6).Value = 1145
$ws.Cells.Item(43, 6).Value = 590
$ws.Cells.Item(45, 6).Value = 3595
$ws.Cells.Item(46, 6).Value = 127
$ws.Cells.Item(48, 6).Value = 1154
$ws.Cells.Item(49, 6).Value = 88
